$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows whose home/away fixture data got reordered (reshuffled rows) ---
# Row 77
$ws.Range("F77").Value = "Goczalkowice Zdroj"
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = "Jelenia Gora"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 1.88
$ws.Range("K77").Value = "07/10/2023 01:43"
$ws.Range("L77").Value = 2.25
$ws.Range("M77").Value = "07/10/2023 12:06"
$ws.Range("N77").Value = 3.52
$ws.Range("O77").Value = "07/10/2023 01:43"
$ws.Range("P77").Value = 3.2
$ws.Range("Q77").Value = "07/10/2023 12:06"
$ws.Range("R77").Value = 3.26
$ws.Range("S77").Value = "07/10/2023 01:43"
$ws.Range("T77").Value = 2.88
$ws.Range("U77").Value = "07/10/2023 12:06"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/goczalkowice-zdroj-karkonosze-jelenia-gora/xQaGA45B/"

# Row 78
$ws.Range("F78").Value = "Starowice Dolne"
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = "Bielsko-Biala"
$ws.Range("I78").Value = 4
$ws.Range("J78").Value = 3.95
$ws.Range("K78").Value = "06/10/2023 00:12"
$ws.Range("L78").Value = 4.67
$ws.Range("M78").Value = "07/10/2023 12:39"
$ws.Range("N78").Value = 3.77
$ws.Range("O78").Value = "06/10/2023 00:12"
$ws.Range("P78").Value = 4.04
$ws.Range("Q78").Value = "07/10/2023 12:39"
$ws.Range("R78").Value = 1.58
$ws.Range("S78").Value = "06/10/2023 00:12"
$ws.Range("T78").Value = 1.55
$ws.Range("U78").Value = "07/10/2023 12:39"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/starowice-dolne-rekord-bielsko-biala/8jeCBpj5/"

# Row 90
$ws.Range("F90").Value = "Goczalkowice Zdroj"
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = "Carina Gubin"
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 1.76
$ws.Range("K90").Value = "20/10/2023 00:12"
$ws.Range("L90").Value = 1.83
$ws.Range("M90").Value = "21/10/2023 12:03"
$ws.Range("N90").Value = 3.41
$ws.Range("O90").Value = "20/10/2023 00:12"
$ws.Range("P90").Value = 3.47
$ws.Range("Q90").Value = "21/10/2023 12:03"
$ws.Range("R90").Value = 3.46
$ws.Range("S90").Value = "20/10/2023 00:12"
$ws.Range("T90").Value = 3.67
$ws.Range("U90").Value = "21/10/2023 12:03"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/goczalkowice-zdroj-carina-gubin/KOLKxNl4/"

# Row 91
$ws.Range("F91").Value = "Starowice Dolne"
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = "Sleza Wroclaw"
$ws.Range("I91").Value = 2
$ws.Range("J91").Value = 3.33
$ws.Range("K91").Value = "20/10/2023 00:12"
$ws.Range("L91").Value = 4.45
$ws.Range("M91").Value = "21/10/2023 12:56"
$ws.Range("N91").Value = 3.75
$ws.Range("O91").Value = "20/10/2023 00:12"
$ws.Range("P91").Value = 4.27
$ws.Range("Q91").Value = "21/10/2023 12:56"
$ws.Range("R91").Value = 1.72
$ws.Range("S91").Value = "20/10/2023 00:12"
$ws.Range("T91").Value = 1.54
$ws.Range("U91").Value = "21/10/2023 12:56"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/starowice-dolne-sleza-wroclaw/tYMGw3Zc/"

# Row 108
$ws.Range("F108").Value = "Gwarek Tarnowskie Gory"
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = "Jelenia Gora"
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 2.6
$ws.Range("K108").Value = "04/11/2023 12:43"
$ws.Range("L108").Value = 2.25
$ws.Range("M108").Value = "04/11/2023 13:10"
$ws.Range("N108").Value = 3.37
$ws.Range("O108").Value = "04/11/2023 12:43"
$ws.Range("P108").Value = 3.42
$ws.Range("Q108").Value = "04/11/2023 13:10"
$ws.Range("R108").Value = 2.31
$ws.Range("S108").Value = "04/11/2023 12:43"
$ws.Range("T108").Value = 2.72
$ws.Range("U108").Value = "04/11/2023 13:10"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gwarek-tarnowskie-gory-karkonosze-jelenia-gora/OxSDpBnr/"

# Row 109
$ws.Range("F109").Value = "Pawlowice"
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = "Rakow II"
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 1.93
$ws.Range("K109").Value = "03/11/2023 02:13"
$ws.Range("L109").Value = 1.88
$ws.Range("M109").Value = "04/11/2023 13:06"
$ws.Range("N109").Value = 3.53
$ws.Range("O109").Value = "03/11/2023 02:13"
$ws.Range("P109").Value = 3.78
$ws.Range("Q109").Value = "04/11/2023 13:06"
$ws.Range("R109").Value = 2.87
$ws.Range("S109").Value = "03/11/2023 02:13"
$ws.Range("T109").Value = 3.25
$ws.Range("U109").Value = "04/11/2023 13:06"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/pniowek-pawlowice-rks-rakow-czestochowa/AgUPs9W0/"

# Row 110
$ws.Range("F110").Value = "Polkowice"
$ws.Range("G110").Value = 3
$ws.Range("H110").Value = "Stilon Gorzow"
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = 1.59
$ws.Range("K110").Value = "03/11/2023 02:13"
$ws.Range("L110").Value = 1.57
$ws.Range("M110").Value = "04/11/2023 13:52"
$ws.Range("N110").Value = 3.85
$ws.Range("O110").Value = "03/11/2023 02:13"
$ws.Range("P110").Value = 4.07
$ws.Range("Q110").Value = "04/11/2023 13:52"
$ws.Range("R110").Value = 3.83
$ws.Range("S110").Value = "03/11/2023 02:13"
$ws.Range("T110").Value = 4.49
$ws.Range("U110").Value = "04/11/2023 13:52"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-stilon-gorzow/xtWHqV1l/"

# Row 111
$ws.Range("F111").Value = "Slask Wroclaw II"
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = "Zielona Gora"
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = 1.53
$ws.Range("K111").Value = "03/11/2023 02:13"
$ws.Range("L111").Value = 1.67
$ws.Range("M111").Value = "04/11/2023 13:09"
$ws.Range("N111").Value = 3.96
$ws.Range("O111").Value = "03/11/2023 02:13"
$ws.Range("P111").Value = 3.99
$ws.Range("Q111").Value = "04/11/2023 13:09"
$ws.Range("R111").Value = 4.1
$ws.Range("S111").Value = "03/11/2023 02:13"
$ws.Range("T111").Value = 3.9
$ws.Range("U111").Value = "04/11/2023 13:09"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/slask-wroclaw-zielona-gora/W2VLrkHf/"

# Row 112
$ws.Range("F112").Value = "Warta Gorzow"
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = "Carina Gubin"
$ws.Range("I112").Value = 1
$ws.Range("J112").Value = 2.41
$ws.Range("K112").Value = "03/11/2023 02:13"
$ws.Range("L112").Value = 2.27
$ws.Range("M112").Value = "04/11/2023 13:42"
$ws.Range("N112").Value = 3.17
$ws.Range("O112").Value = "03/11/2023 02:13"
$ws.Range("P112").Value = 3.16
$ws.Range("Q112").Value = "04/11/2023 13:21"
$ws.Range("R112").Value = 2.4
$ws.Range("S112").Value = "03/11/2023 02:13"
$ws.Range("T112").Value = 2.87
$ws.Range("U112").Value = "04/11/2023 13:42"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/warta-gorzow-carina-gubin/4lFu3J29/"

# Row 116
$ws.Range("F116").Value = "Zielona Gora"
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = "Pawlowice"
$ws.Range("I116").Value = 1
$ws.Range("J116").Value = 2.06
$ws.Range("K116").Value = "11/11/2023 01:13"
$ws.Range("L116").Value = 1.84
$ws.Range("M116").Value = "11/11/2023 08:02"
$ws.Range("N116").Value = 3.48
$ws.Range("O116").Value = "11/11/2023 01:13"
$ws.Range("P116").Value = 3.59
$ws.Range("Q116").Value = "11/11/2023 11:04"
$ws.Range("R116").Value = 2.87
$ws.Range("S116").Value = "11/11/2023 01:13"
$ws.Range("T116").Value = 3.47
$ws.Range("U116").Value = "11/11/2023 08:02"
$ws.Range("V116").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/zielona-gora-pniowek-pawlowice/UepRum1D/"

# Row 117
$ws.Range("F117").Value = "Unia Turza Slaska"
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = "Bytom Odrzanski"
$ws.Range("I117").Value = 2
$ws.Range("J117").Value = 1.93
$ws.Range("K117").Value = "11/11/2023 02:12"
$ws.Range("L117").Value = 2.02
$ws.Range("M117").Value = "11/11/2023 12:23"
$ws.Range("N117").Value = 3.61
$ws.Range("O117").Value = "11/11/2023 02:12"
$ws.Range("P117").Value = 3.41
$ws.Range("Q117").Value = "11/11/2023 12:23"
$ws.Range("R117").Value = 3.05
$ws.Range("S117").Value = "11/11/2023 02:12"
$ws.Range("T117").Value = 3.15
$ws.Range("U117").Value = "11/11/2023 12:23"
$ws.Range("V117").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/unia-turza-slaska-bytom-odrzanski/l0Hid8of/"

# Row 119
$ws.Range("F119").Value = "Bielsko-Biala"
$ws.Range("G119").Value = 3
$ws.Range("H119").Value = "Gwarek Tarnowskie Gory"
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = 1.25
$ws.Range("K119").Value = "11/11/2023 01:13"
$ws.Range("L119").Value = 1.37
$ws.Range("M119").Value = "11/11/2023 12:44"
$ws.Range("N119").Value = 5.43
$ws.Range("O119").Value = "11/11/2023 01:13"
$ws.Range("P119").Value = 5.47
$ws.Range("Q119").Value = "11/11/2023 12:44"
$ws.Range("R119").Value = 7.09
$ws.Range("S119").Value = "11/11/2023 01:13"
$ws.Range("T119").Value = 5.17
$ws.Range("U119").Value = "11/11/2023 12:44"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/rekord-bielsko-biala-gwarek-tarnowskie-gory/vRRH2A8Q/"

# Row 120
$ws.Range("F120").Value = "Carina Gubin"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "Gornik Zabrze II"
$ws.Range("I120").Value = 2
$ws.Range("J120").Value = 2.12
$ws.Range("K120").Value = "11/11/2023 02:12"
$ws.Range("L120").Value = 2.31
$ws.Range("M120").Value = "11/11/2023 12:52"
$ws.Range("N120").Value = 3.58
$ws.Range("O120").Value = "11/11/2023 02:12"
$ws.Range("P120").Value = 3.52
$ws.Range("Q120").Value = "11/11/2023 12:52"
$ws.Range("R120").Value = 2.7
$ws.Range("S120").Value = "11/11/2023 02:12"
$ws.Range("T120").Value = 2.58
$ws.Range("U120").Value = "11/11/2023 12:52"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/carina-gubin-gornik-zabrze/tzEqbUGs/"

# --- Append new row 131 (new match result) ---
$ws.Range("A130").Copy()
$ws.Range("A131").PasteSpecial(-4122)
$ws.Range("E130").Copy()
$ws.Range("E131").PasteSpecial(-4122)

$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "poland"
$ws.Range("C131").Value = "iii-liga-group-iii"
$ws.Range("D131").Value = "2023-2024"
$ws.Range("E131").Value = 45254.70833333334
$ws.Range("F131").Value = "Kluczbork"
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = "Zielona Gora"
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1.87
$ws.Range("K131").Value = "24/11/2023 06:12"
$ws.Range("L131").Value = 2.08
$ws.Range("M131").Value = "24/11/2023 10:19"
$ws.Range("N131").Value = 3.63
$ws.Range("O131").Value = "24/11/2023 06:12"
$ws.Range("P131").Value = 3.71
$ws.Range("Q131").Value = "24/11/2023 15:03"
$ws.Range("R131").Value = 3.3
$ws.Range("S131").Value = "24/11/2023 06:12"
$ws.Range("T131").Value = 2.81
$ws.Range("U131").Value = "24/11/2023 16:26"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/kluczbork-zielona-gora/Igpw85p6/"
